# Update probability matrix values on Sheet1 (Ohio_B team-specific matrix)
# Reflects recomputed simulation statistics after adding more games,
# speeding up simulate-game logic, and drafting optimization logic.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

    $ws.Range("B2").Value = 0.1606334841628959
    $ws.Range("C2").Value = 0.6153846153846154
    $ws.Range("J2").Value = 0.006787330316742082
    $ws.Range("P2").Value = 0.1266968325791855
    $ws.Range("S2").Value = 0.09049773755656108
    $ws.Range("B3").Value = 0.003597122302158274
    $ws.Range("C3").Value = 0.02158273381294964
    $ws.Range("J3").Value = 0.02517985611510791
    $ws.Range("P3").Value = 0.7266187050359713
    $ws.Range("S3").Value = 0.223021582733813
    $ws.Range("J4").Value = 0.03125
    $ws.Range("P4").Value = 0.640625
    $ws.Range("S4").Value = 0.328125
    $ws.Range("P5").Value = 1
    $ws.Range("B6").Value = 0.06191950464396285
    $ws.Range("D6").Value = 0.01547987616099071
    $ws.Range("E6").Value = 0.003095975232198143
    $ws.Range("F6").Value = 0.06191950464396285
    $ws.Range("J6").Value = 0.238390092879257
    $ws.Range("O6").Value = 0.01238390092879257
    $ws.Range("Q6").Value = 0.173374613003096
    $ws.Range("R6").Value = 0.06811145510835913
    $ws.Range("S6").Value = 0.3653250773993808
    $ws.Range("B7").Value = 0.1295681063122923
    $ws.Range("D7").Value = 0.02325581395348837
    $ws.Range("F7").Value = 0.05315614617940199
    $ws.Range("J7").Value = 0.1561461794019934
    $ws.Range("O7").Value = 0.009966777408637873
    $ws.Range("Q7").Value = 0.1661129568106312
    $ws.Range("R7").Value = 0.09966777408637874
    $ws.Range("S7").Value = 0.3621262458471761
    $ws.Range("B8").Value = 0.1051873198847262
    $ws.Range("D8").Value = 0.02161383285302594
    $ws.Range("F8").Value = 0.0590778097982709
    $ws.Range("J8").Value = 0.1210374639769452
    $ws.Range("O8").Value = 0.01296829971181556
    $ws.Range("Q8").Value = 0.1858789625360231
    $ws.Range("R8").Value = 0.09221902017291066
    $ws.Range("S8").Value = 0.4020172910662824
    $ws.Range("B9").Value = 0.08888888888888889
    $ws.Range("D9").Value = 0.0126984126984127
    $ws.Range("F9").Value = 0.07936507936507936
    $ws.Range("J9").Value = 0.1047619047619048
    $ws.Range("O9").Value = 0.0380952380952381
    $ws.Range("Q9").Value = 0.2
    $ws.Range("R9").Value = 0.09206349206349207
    $ws.Range("S9").Value = 0.3841269841269842
    $ws.Range("B10").Value = 0.1103260869565217
    $ws.Range("D10").Value = 0.01956521739130435
    $ws.Range("F10").Value = 0.07065217391304347
    $ws.Range("J10").Value = 0.1309782608695652
    $ws.Range("O10").Value = 0.01032608695652174
    $ws.Range("Q10").Value = 0.2097826086956522
    $ws.Range("R10").Value = 0.08043478260869565
    $ws.Range("S10").Value = 0.3679347826086957
    $ws.Range("G11").Value = 0.1461716937354988
    $ws.Range("J11").Value = 0.07888631090487239
    $ws.Range("K11").Value = 0.1879350348027842
    $ws.Range("L11").Value = 0.5707656612529002
    $ws.Range("S11").Value = 0.01624129930394431
    $ws.Range("G12").Value = 0.7831325301204819
    $ws.Range("J12").Value = 0.1807228915662651
    $ws.Range("K12").Value = 0.008032128514056224
    $ws.Range("L12").Value = 0.02008032128514056
    $ws.Range("S12").Value = 0.008032128514056224
    $ws.Range("G13").Value = 0.847457627118644
    $ws.Range("J13").Value = 0.1186440677966102
    $ws.Range("S13").Value = 0.03389830508474576
    $ws.Range("F15").Value = 0.01257861635220126
    $ws.Range("H15").Value = 0.1918238993710692
    $ws.Range("I15").Value = 0.1069182389937107
    $ws.Range("J15").Value = 0.3176100628930817
    $ws.Range("K15").Value = 0.05031446540880503
    $ws.Range("M15").Value = 0.01572327044025157
    $ws.Range("O15").Value = 0.0660377358490566
    $ws.Range("S15").Value = 0.2389937106918239
    $ws.Range("F16").Value = 0.02068965517241379
    $ws.Range("H16").Value = 0.2172413793103448
    $ws.Range("I16").Value = 0.06206896551724138
    $ws.Range("J16").Value = 0.3896551724137931
    $ws.Range("K16").Value = 0.1137931034482759
    $ws.Range("M16").Value = 0.02068965517241379
    $ws.Range("O16").Value = 0.04827586206896552
    $ws.Range("S16").Value = 0.1275862068965517
    $ws.Range("F17").Value = 0.0117820324005891
    $ws.Range("H17").Value = 0.1899852724594993
    $ws.Range("I17").Value = 0.101620029455081
    $ws.Range("J17").Value = 0.3652430044182621
    $ws.Range("K17").Value = 0.101620029455081
    $ws.Range("M17").Value = 0.007363770250368188
    $ws.Range("N17").Value = 0.002945508100147275
    $ws.Range("O17").Value = 0.07511045655375552
    $ws.Range("S17").Value = 0.1443298969072165
    $ws.Range("F18").Value = 0.02033898305084746
    $ws.Range("H18").Value = 0.1627118644067796
    $ws.Range("I18").Value = 0.0711864406779661
    $ws.Range("J18").Value = 0.4406779661016949
    $ws.Range("K18").Value = 0.1220338983050848
    $ws.Range("M18").Value = 0.01016949152542373
    $ws.Range("O18").Value = 0.05084745762711865
    $ws.Range("S18").Value = 0.1220338983050848
    $ws.Range("F19").Value = 0.01183431952662722
    $ws.Range("H19").Value = 0.2157073695535234
    $ws.Range("I19").Value = 0.09306078536847767
    $ws.Range("J19").Value = 0.3722431414739107
    $ws.Range("K19").Value = 0.1032813340505648
    $ws.Range("M19").Value = 0.02205486820871436
    $ws.Range("N19").Value = 0.0005379236148466917
    $ws.Range("O19").Value = 0.07692307692307693
    $ws.Range("S19").Value = 0.1043571812802582

$wb.Save()
